# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (interest count) figures in the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13074
$ws1.Range("F8").Value = 29
$ws1.Range("F10").Value = 13039
$ws1.Range("F11").Value = 298
$ws1.Range("F15").Value = 212

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13074
$ws4.Range("F9").Value = 29
$ws4.Range("F11").Value = 13039
$ws4.Range("F12").Value = 298
$ws4.Range("F16").Value = 212
